{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph ending with the \"third model\" colorization sentence --\n// the anchor point after which the new \"After T59v5\" section is inserted.\nconst anchorMarker = \"out-of-sequence.\";\nlet anchor = null;\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(anchorMarker) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\nconst newBodyXml = `<w:p><w:r><w:br w:type=\"page\"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>After T59v5</w:t></w:r></w:p><w:p><w:r><w:t>Hierarchical</w:t></w:r><w:r><w:t xml:space=\"preserve\"> structuration of</w:t></w:r><w:r><w:t xml:space=\"preserve\"> T59 </w:t></w:r><w:r><w:t>source code:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Vocab token</w:t></w:r><w:r><w:t>/</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>IToken</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>/L0Token</w:t></w:r><w:r><w:t>: sequence returned by ANTLR lexical analysis, including white space on a separate channel. The full sequence covers 100% of source code.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L1Token</w:t></w:r><w:r><w:t xml:space=\"preserve\">: Ignore WS, group successive </w:t></w:r><w:r><w:t>INVALID_CHAR</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">into </w:t></w:r><w:r><w:t>L1InvalidToken</w:t></w:r><w:r><w:t xml:space=\"preserve\">, transform all </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>I_xx</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> into L1Instruction with </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>TIKey</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> attribute, has a </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>SyntaxCategory</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>. Note that L1A4 token (build during L2 pass) groups two consecutive Vocab.D2 when they should be interpreted as an address. Sequence of L1Tokens covers source code ignoring white space.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L2Statement</w:t></w:r><w:r><w:t xml:space=\"preserve\">: group of L1Tokens, including L2InvalidStatement, and two </w:t></w:r><w:r><w:t>L2ActualInstruction</w:t></w:r><w:r><w:t xml:space=\"preserve\">: </w:t></w:r><w:r><w:t>L2Instruction</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and L2Number with </w:t></w:r><w:r><w:t>A</w:t></w:r><w:r><w:t>ddress</w:t></w:r><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>OpCodes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> and Problem flag (detected during L3 encoding, such as invalid address or duplicate label). L2Tag has also an Address and a Problem flag.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>T59Program</w:t></w:r><w:r><w:t>: group of L2Statements, valid addresses and a list of errors detected during L2encoding. Maybe a T59Error class should be used instead of a string, and replace L2Statement Problem flag by a reference to T59Error.</w:t></w:r></w:p><w:p><w:r><w:t>When colorizing source code preserving its integrity, work on the sequence of L0Tokens: WS ignored, non-WS: L0Token-&gt;L1Token.SyntaxCategory and L0Token-&gt;L1Token-&gt;L2Statement is L2InvalidStatement determines color, and L0Token-&gt;L1Token-&gt;L2Statement.Problem may add a visual flag.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">When reformatting, start on list of L2Statements, print </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>OpCodes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> and reformatted</w:t></w:r><w:r><w:t>/colored</w:t></w:r><w:r><w:t xml:space=\"preserve\"> L1Statements</w:t></w:r><w:r><w:t xml:space=\"preserve\"> (canonical mnemonics for instance), standardized white space.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>To colorize source code, I need a quick way to move up in the hierarchy: L0Token -&gt; L1Token -&gt; L2Statement.</w:t></w:r></w:p>`;\n\nconst ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newBodyXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst insertRange = anchor.getRange(Word.RangeLocation.end);\ninsertRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph ending with the \"third model\" colorization sentence --\n# the anchor point after which the new \"After T59v5\" section is inserted.\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*out-of-sequence*\") {\n        $anchor = $p\n    }\n}\nif ($null -eq $anchor) {\n    throw \"Anchor paragraph not found\"\n}\n\n$insPoint = $anchor.Range.End\n$insertionRange = $d.Range($insPoint, $insPoint)\n\n$newBodyXml = '<w:p><w:r><w:br w:type=\"page\"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>After T59v5</w:t></w:r></w:p><w:p><w:r><w:t>Hierarchical</w:t></w:r><w:r><w:t xml:space=\"preserve\"> structuration of</w:t></w:r><w:r><w:t xml:space=\"preserve\"> T59 </w:t></w:r><w:r><w:t>source code:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t>Vocab token</w:t></w:r><w:r><w:t>/</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>IToken</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>/L0Token</w:t></w:r><w:r><w:t>: sequence returned by ANTLR lexical analysis, including white space on a separate channel. The full sequence covers 100% of source code.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L1Token</w:t></w:r><w:r><w:t xml:space=\"preserve\">: Ignore WS, group successive </w:t></w:r><w:r><w:t>INVALID_CHAR</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">into </w:t></w:r><w:r><w:t>L1InvalidToken</w:t></w:r><w:r><w:t xml:space=\"preserve\">, transform all </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>I_xx</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> into L1Instruction with </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>TIKey</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> attribute, has a </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>SyntaxCategory</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>. Note that L1A4 token (build during L2 pass) groups two consecutive Vocab.D2 when they should be interpreted as an address. Sequence of L1Tokens covers source code ignoring white space.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L2Statement</w:t></w:r><w:r><w:t xml:space=\"preserve\">: group of L1Tokens, including L2InvalidStatement, and two </w:t></w:r><w:r><w:t>L2ActualInstruction</w:t></w:r><w:r><w:t xml:space=\"preserve\">: </w:t></w:r><w:r><w:t>L2Instruction</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and L2Number with </w:t></w:r><w:r><w:t>A</w:t></w:r><w:r><w:t>ddress</w:t></w:r><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>OpCodes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> and Problem flag (detected during L3 encoding, such as invalid address or duplicate label). L2Tag has also an Address and a Problem flag.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>T59Program</w:t></w:r><w:r><w:t>: group of L2Statements, valid addresses and a list of errors detected during L2encoding. Maybe a T59Error class should be used instead of a string, and replace L2Statement Problem flag by a reference to T59Error.</w:t></w:r></w:p><w:p><w:r><w:t>When colorizing source code preserving its integrity, work on the sequence of L0Tokens: WS ignored, non-WS: L0Token-&gt;L1Token.SyntaxCategory and L0Token-&gt;L1Token-&gt;L2Statement is L2InvalidStatement determines color, and L0Token-&gt;L1Token-&gt;L2Statement.Problem may add a visual flag.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">When reformatting, start on list of L2Statements, print </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>OpCodes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> and reformatted</w:t></w:r><w:r><w:t>/colored</w:t></w:r><w:r><w:t xml:space=\"preserve\"> L1Statements</w:t></w:r><w:r><w:t xml:space=\"preserve\"> (canonical mnemonics for instance), standardized white space.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>To colorize source code, I need a quick way to move up in the hierarchy: L0Token -&gt; L1Token -&gt; L2Statement.</w:t></w:r></w:p><w:p/>'\n\n$ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$insertionRange.InsertXML($ooxml)\n"}
